$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`r`n✅ Dólar paralelo: 68`r`n`r`nBinance`r`n✅ 1000 Bs = 5.14 = 20498.77 pesos`r`n✅ 20498.77 pesos = 5.12 = 968.08 Bs`r`n`r`nPromedio competencia`r`n✅ Tasa pesos: 20`r`n✅ Tasa Bs: 20`r`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 194.5
$ws2.Range("O10").Value = 3987.01
$ws2.Range("N12").Value = 4002
$ws2.Range("O12").Value = 189
